# Apply crypto price/volume updates (and OKB/Polygon, Hedera/VeChain/WrappedliquidstakedEther2.0 row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.844.39"
$ws.Range("E2").Value = "  +0.99%  "

# Row 3
$ws.Range("D3").Value = "1.706.37"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.96"
$ws.Range("E5").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4017"
$ws.Range("E7").Value = "  +3.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4042"
$ws.Range("E8").Value = "  +0.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.003"
$ws.Range("E9").Value = "  -0.26%  "

# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.77"
$ws.Range("E10").Value = "  +1.84%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.469"
$ws.Range("E11").Value = "  -1.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08794"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.36"
$ws.Range("E13").Value = "  +6.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.500"
$ws.Range("E14").Value = "  -1.38%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.000"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001341"
$ws.Range("E16").Value = "  +0.31%  "

# Row 17
$ws.Range("D17").Value = "1.610.16"
$ws.Range("E17").Value = "  -4.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.42"
$ws.Range("E18").Value = "  -2.62%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07178"
$ws.Range("E19").Value = "  +1.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.92"
$ws.Range("E20").Value = "  +6.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.286"
$ws.Range("E21").Value = "  +0.70%  "

# Row 22
$ws.Range("E22").Value = "  -0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.44"
$ws.Range("E23").Value = "  +1.98%  "

# Row 24
$ws.Range("D24").Value = "24.845.15"
$ws.Range("E24").Value = "  +0.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.339"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.889"
$ws.Range("E26").Value = "  -3.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.411"
$ws.Range("E27").Value = "  +22.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.09"
$ws.Range("E28").Value = "  +2.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.70"
$ws.Range("E29").Value = "  +0.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "143.86"
$ws.Range("E30").Value = "  +5.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.321"
$ws.Range("E31").Value = "  -2.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.290"
$ws.Range("E32").Value = "  +15.81%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08681"
$ws.Range("E33").Value = "  -1.01%  "

# Row 34
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03183"
$ws.Range("E34").Value = "  +10.19%  "

# Row 35
$ws.Range("B35").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C35").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D35").Value = "1.788.64"
$ws.Range("E35").Value = "  -4.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.211"
$ws.Range("E36").Value = "  -2.78%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.028"
$ws.Range("E37").Value = "  -0.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2860"
$ws.Range("E38").Value = "  +5.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8409"
$ws.Range("E39").Value = "  +8.51%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.82"
$ws.Range("E40").Value = "  +0.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09432"
$ws.Range("E41").Value = "  +3.75%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.23"
$ws.Range("E42").Value = "  +1.58%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.481"
$ws.Range("E43").Value = "  +1.81%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.44"
$ws.Range("E44").Value = "  +5.87%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.734"
$ws.Range("E45").Value = "  +6.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7417"
$ws.Range("E46").Value = "  +4.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.224"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.365"
$ws.Range("E48").Value = "  +1.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.43"
$ws.Range("E50").Value = "  +1.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08390"
$ws.Range("E51").Value = "  +5.25%  "

